$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 0.9584923333333334
$ws.Range("H2").Value = 2.875477
$ws.Range("I2").Value = 0.2532195598902293
$ws.Range("J2").Value = 0.2532195598902293
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 14.29506233333333
$ws.Range("N2").Value = 42.885187
$ws.Range("O2").Value = 0.2033714702773193
$ws.Range("P2").Value = 0.2033714702773193
$ws.Range("Q2").Value = 13.70170765102211
$ws.Range("R2").Value = 123.315368859199
$ws.Range("S2").Value = 0.05149763419785165
$ws.Range("T2").Value = 0.05149763419785163

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 0.9584923333333334
$ws.Range("H3").Value = 2.875477
$ws.Range("I3").Value = 0.2532195598902293
$ws.Range("J3").Value = 0.2532195598902293
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 21.21615033333333
$ws.Range("N3").Value = 63.648451
$ws.Range("O3").Value = 0.3018356678902651
$ws.Range("P3").Value = 0.3018356678902651
$ws.Range("Q3").Value = 20.33551743734745
$ws.Range("R3").Value = 183.019656936127
$ws.Range("S3").Value = 0.07643069498234634
$ws.Range("T3").Value = 0.07643069498234632

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 0.9584923333333334
$ws.Range("H4").Value = 2.875477
$ws.Range("I4").Value = 0.2532195598902293
$ws.Range("J4").Value = 0.2532195598902293
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 14.99890066666667
$ws.Range("N4").Value = 44.996702
$ws.Range("O4").Value = 0.2133847625141612
$ws.Range("P4").Value = 0.2133847625141612
$ws.Range("Q4").Value = 14.37633129742822
$ws.Range("R4").Value = 129.386981676854
$ws.Range("S4").Value = 0.05403319565111699
$ws.Range("T4").Value = 0.05403319565111699

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 0.9584923333333334
$ws.Range("H5").Value = 2.875477
$ws.Range("I5").Value = 0.2532195598902293
$ws.Range("J5").Value = 0.2532195598902293
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 16.03541933333333
$ws.Range("N5").Value = 48.106258
$ws.Range("O5").Value = 0.228130995884431
$ws.Range("P5").Value = 0.2281309958844309
$ws.Range("Q5").Value = 15.36982649278511
$ws.Range("R5").Value = 138.328438435066
$ws.Range("S5").Value = 0.05776723037517532
$ws.Range("T5").Value = 0.05776723037517531

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 0.9584923333333334
$ws.Range("H6").Value = 2.875477
$ws.Range("I6").Value = 0.2532195598902293
$ws.Range("J6").Value = 0.2532195598902293
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 3.744869
$ws.Range("N6").Value = 11.234607
$ws.Range("O6").Value = 0.05327710343382351
$ws.Range("P6").Value = 0.0532771034338235
$ws.Range("Q6").Value = 3.589428225837667
$ws.Range("R6").Value = 32.30485403253901
$ws.Range("S6").Value = 0.01349080468373901
$ws.Range("T6").Value = 0.01349080468373901

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 2.82673
$ws.Range("H7").Value = 8.48019
$ws.Range("I7").Value = 0.7467804401097707
$ws.Range("J7").Value = 0.7467804401097707
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 14.29506233333333
$ws.Range("N7").Value = 42.885187
$ws.Range("O7").Value = 0.2033714702773193
$ws.Range("P7").Value = 0.2033714702773193
$ws.Range("Q7").Value = 40.40828154950334
$ws.Range("R7").Value = 363.67453394553
$ws.Range("S7").Value = 0.1518738360794677
$ws.Range("T7").Value = 0.1518738360794676

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 2.82673
$ws.Range("H8").Value = 8.48019
$ws.Range("I8").Value = 0.7467804401097707
$ws.Range("J8").Value = 0.7467804401097707
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 21.21615033333333
$ws.Range("N8").Value = 63.648451
$ws.Range("O8").Value = 0.3018356678902651
$ws.Range("P8").Value = 0.3018356678902651
$ws.Range("Q8").Value = 59.97232863174334
$ws.Range("R8").Value = 539.7509576856901
$ws.Range("S8").Value = 0.2254049729079188
$ws.Range("T8").Value = 0.2254049729079187

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 2.82673
$ws.Range("H9").Value = 8.48019
$ws.Range("I9").Value = 0.7467804401097707
$ws.Range("J9").Value = 0.7467804401097707
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 14.99890066666667
$ws.Range("N9").Value = 44.996702
$ws.Range("O9").Value = 0.2133847625141612
$ws.Range("P9").Value = 0.2133847625141612
$ws.Range("Q9").Value = 42.39784248148666
$ws.Range("R9").Value = 381.58058233338
$ws.Range("S9").Value = 0.1593515668630442
$ws.Range("T9").Value = 0.1593515668630442

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 2.82673
$ws.Range("H10").Value = 8.48019
$ws.Range("I10").Value = 0.7467804401097707
$ws.Range("J10").Value = 0.7467804401097707
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 16.03541933333333
$ws.Range("N10").Value = 48.106258
$ws.Range("O10").Value = 0.228130995884431
$ws.Range("P10").Value = 0.2281309958844309
$ws.Range("Q10").Value = 45.32780089211333
$ws.Range("R10").Value = 407.95020802902
$ws.Range("S10").Value = 0.1703637655092557
$ws.Range("T10").Value = 0.1703637655092556

# Row 11
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 2.82673
$ws.Range("H11").Value = 8.48019
$ws.Range("I11").Value = 0.7467804401097707
$ws.Range("J11").Value = 0.7467804401097707
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 3.744869
$ws.Range("N11").Value = 11.234607
$ws.Range("O11").Value = 0.05327710343382351
$ws.Range("P11").Value = 0.0532771034338235
$ws.Range("Q11").Value = 10.58573354837
$ws.Range("R11").Value = 95.27160193533001
$ws.Range("S11").Value = 0.0397862987500845
$ws.Range("T11").Value = 0.03978629875008449

Write-Host "Applied NATMI recalculated values (Dr Hou advice)"